# Commit: "urgent messages and new message"
#
# On the "Suite" sheet, flip the Runmode ("Y"/"N") flags for two test
# cases:
#   - Web_QUICK_MESSAGES (row 4)  : N -> Y
#   - NEW_MESSAGE        (row 17) : Y -> N
# and leave the cursor/selection on E10, matching where the author's
# Excel session ended up after making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Web_QUICK_MESSAGES row: Runmode N -> Y
$ws.Range("C4").Value = "Y"

# NEW_MESSAGE row: Runmode Y -> N
$ws.Range("C17").Value = "N"

# Leave the selection where the author left it.
$ws.Range("E10").Select() | Out-Null
